$wb = $excel.ActiveWorkbook

# Overview sheet: zh-cn / de-de status columns and the overall latest handoff date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-22 06:38:16"

# zh-cn sheet: Status + Latest Handoff Datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-22 06:38:13"

# de-de sheet: Status + Latest Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-22 06:38:16"
